# Rename the existing sheet (the big incident grid) to "Table"
$wb = $excel.ActiveWorkbook
$origWs = $wb.Worksheets.Item(1)
$origWs.Name = "Table"

# Add a brand-new sheet; since "Table" took the old "Sheet1" name, Add()
# will place this new sheet first and auto-name it "Sheet1"
$newWs = $wb.Worksheets.Add()

# Re-fetch stable references to both sheets by name (object references
# captured before Add() can report stale index/name info)
$sheet1 = $wb.Worksheets.Item("Sheet1")
$tableWs = $wb.Worksheets.Item("Table")

# Populate the new "Sheet1" with the Month / Year summary table
$sheet1.Range("A1").Value = "Month"
$sheet1.Range("B1").Value = "Year"

$sheet1.Cells.Item(2, 1).Value = 11
$sheet1.Cells.Item(2, 2).Value = 2019

$sheet1.Cells.Item(3, 1).Value = 4
$sheet1.Cells.Item(3, 2).Value = 2020

$sheet1.Cells.Item(4, 1).Value = 8
$sheet1.Cells.Item(4, 2).Value = 2020

$sheet1.Cells.Item(5, 1).Value = 5
$sheet1.Cells.Item(5, 2).Value = 2021

$sheet1.Cells.Item(6, 1).Value = 8
$sheet1.Cells.Item(6, 2).Value = 2022

$sheet1.Cells.Item(7, 1).Value = 8
$sheet1.Cells.Item(7, 2).Value = 2023

# Restore the selection on the "Table" sheet and make "Sheet1" the
# active / selected tab, matching the target workbook state
$tableWs.Range("C21").Select()
$sheet1.Range("D5").Select()
